$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.334.90'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.036.39'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.71%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.91'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.37%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.392'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0806'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.17'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.340.84'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.848'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.95'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.040.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.278.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +4.38%  '
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.18%  '
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.91%  '
$ws.Range('E30').Value = '  +3.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.121'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0680'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.80%  '
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.04%  '
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('E36').Value = '  +6.02%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('E38').Value = '  +2.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('E40').Value = '  +2.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0979'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('E42').Value = '  +3.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.62'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.22%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.388.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.06'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.232.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.04%  '
